# Generate Report for Handoff
#
# a76b8f79-... now reports first (row 2) on every sheet, 1daa9765-...
# drops to row 3 and gets a fresh (not-yet-in-sync) handoff status, an
# updated handoff datetime, and a new "stale handback" error detail.
#
# Only cells whose value actually changes are written so that pre-existing
# blank cells (e.g. the "Dependency From" / unused columns) are left
# exactly as they were instead of being collapsed/removed.

$wb = $excel.ActiveWorkbook

$u1 = "1daa9765-d5e8-4f9f-a8cd-e589105224c0"
$u2 = "a76b8f79-3bfc-467c-80fd-cf3c63741b28"

$urlBase   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f695303534129e140d551b0a0e6aa92cbc97e0ab/e2e/"
$urlZhBase = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e32aab4a8f7b02195bd627b97e80fa0b5a9d057c/e2e/"
$urlDeBase = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/ad39d4e5acdebab8a5e737b57da49af4f84f327c/e2e/"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f695303534129e140d551b0a0e6aa92cbc97e0ab/e2e/$u1.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4b8d3d1048b89dcfd7319bc12dc0f5596ecddd7d/e2e/$u1.md."

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "$u2.md"
$ov.Range("B2").Value = "e2e\$u2.md"

$ov.Range("A3").Value = "$u1.md"
$ov.Range("B3").Value = "e2e\$u1.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-29 10:49:26"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "$urlBase$u1.md", "", "", "e2e\$u2.md")
$ov.Hyperlinks.Add($ov.Range("B3"), "$urlBase$u2.md", "", "", "e2e\$u1.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "$u2.md"
$zh.Range("G2").Value = "$u2.31bd3791c66b9e7ee668cc36dde2633eab440f78.zh-cn.xlf"
$zh.Range("I2").Value = "$u2.md"
$zh.Range("J2").Value = "$u2.31bd3791c66b9e7ee668cc36dde2633eab440f78.zh-cn.xlf"

$zh.Range("A3").Value = "$u1.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "$u1.dc57d2b065390c90d27816386da34daf2d8b263f.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-29 10:49:22"
$zh.Range("I3").Value = "$u1.md"
$zh.Range("J3").Value = "$u1.dc57d2b065390c90d27816386da34daf2d8b263f.zh-cn.xlf"
$zh.Range("P3").Value = $errorDetail

$zh.Columns.Item(16).ColumnWidth = 40

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "$urlBase$u1.md", "", "", "$u2.md")
$zh.Hyperlinks.Add($zh.Range("I2"), "$urlZhBase$u1.md", "", "", "$u2.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "$urlBase$u2.md", "", "", "$u1.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "$urlZhBase$u2.md", "", "", "$u1.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "$u2.md"
$de.Range("G2").Value = "$u2.31bd3791c66b9e7ee668cc36dde2633eab440f78.de-de.xlf"
$de.Range("I2").Value = "$u2.md"
$de.Range("J2").Value = "$u2.31bd3791c66b9e7ee668cc36dde2633eab440f78.de-de.xlf"

$de.Range("A3").Value = "$u1.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "$u1.dc57d2b065390c90d27816386da34daf2d8b263f.de-de.xlf"
$de.Range("H3").Value = "2016-08-29 10:49:26"
$de.Range("I3").Value = "$u1.md"
$de.Range("J3").Value = "$u1.dc57d2b065390c90d27816386da34daf2d8b263f.de-de.xlf"
$de.Range("P3").Value = $errorDetail

$de.Columns.Item(16).ColumnWidth = 40

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "$urlBase$u1.md", "", "", "$u2.md")
$de.Hyperlinks.Add($de.Range("I2"), "$urlDeBase$u1.md", "", "", "$u2.md")
$de.Hyperlinks.Add($de.Range("A3"), "$urlBase$u2.md", "", "", "$u1.md")
$de.Hyperlinks.Add($de.Range("I3"), "$urlDeBase$u2.md", "", "", "$u1.md")

Write-Host "Done applying handoff report changes."
